# Auto-generated Excel COM-interop edit script
# Applies the "Updated cryptos list" diff: refreshed prices/volume percentages
# and a block of coin rows (41-46) that got reshuffled/re-ranked.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.605.20'
$ws.Range('E2').Value = '  +5.01%  '
$ws.Range('D3').Value = '3.500.75'
$ws.Range('E3').Value = '  +3.14%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '589.87'
$ws.Range('E5').Value = '  +4.07%  '
$ws.Range('D6').Value = '169.32'
$ws.Range('E6').Value = '  +8.41%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.502.08'
$ws.Range('E8').Value = '  +3.14%  '
$ws.Range('D9').Value = '''0.580'
$ws.Range('E9').Value = '  +2.48%  '
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('E11').Value = '  +4.80%  '
$ws.Range('E12').Value = '  +2.76%  '
$ws.Range('D13').Value = '4.105.30'
$ws.Range('E13').Value = '  +3.16%  '
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('E15').Value = '  +4.09%  '
$ws.Range('D16').Value = '66.567.22'
$ws.Range('E16').Value = '  +4.84%  '
$ws.Range('E17').Value = '  +3.60%  '
$ws.Range('D18').Value = '3.501.12'
$ws.Range('E18').Value = '  +2.87%  '
$ws.Range('E19').Value = '  +3.20%  '
$ws.Range('D20').Value = '13.98'
$ws.Range('E20').Value = '  +3.52%  '
$ws.Range('D21').Value = '''388.20'
$ws.Range('E21').Value = '  +1.13%  '
$ws.Range('E22').Value = '  +3.27%  '
$ws.Range('D23').Value = '72.96'
$ws.Range('E23').Value = '  +2.75%  '
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('E25').Value = '  +2.39%  '
$ws.Range('E26').Value = '  +8.28%  '
$ws.Range('D27').Value = '10.11'
$ws.Range('E27').Value = '  +4.49%  '
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  +4.87%  '
$ws.Range('E31').Value = '  +5.56%  '
$ws.Range('E32').Value = '  +4.17%  '
$ws.Range('D33').Value = '23.43'
$ws.Range('E33').Value = '  +2.77%  '
$ws.Range('E34').Value = '  +7.29%  '
$ws.Range('E36').Value = '  +4.07%  '
$ws.Range('D37').Value = '161.13'
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('D38').Value = '0.902'
$ws.Range('E38').Value = '  +7.75%  '
$ws.Range('E39').Value = '  +6.11%  '
$ws.Range('D40').Value = '0.0746'
$ws.Range('E40').Value = '  +4.03%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '4.63'
$ws.Range('E41').Value = '  +6.63%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '''6.70'
$ws.Range('E42').Value = '  +4.60%  '
$ws.Range('D43').Value = '26.37'
$ws.Range('E43').Value = '  +2.19%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '43.37'
$ws.Range('E44').Value = '  +1.22%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.794.66'
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '26.72'
$ws.Range('E46').Value = '  +4.58%  '
$ws.Range('E47').Value = '  +3.61%  '
$ws.Range('E48').Value = '  +9.41%  '
$ws.Range('D49').Value = '353.36'
$ws.Range('E49').Value = '  +8.37%  '
$ws.Range('E50').Value = '  +6.38%  '
$ws.Range('D51').Value = '33.12'
$ws.Range('E51').Value = '  +9.92%  '
